# Auto-generated: update market-derived profit columns (H-N) per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, "H").Value = 144.85
$ws.Cells.Item(33, "I").Value = 152.76471
$ws.Cells.Item(33, "K").Value = 152.76471
$ws.Cells.Item(33, "M").Value = 76.23528999999999

# Row 92 (Leve Item ID 19901)
$ws.Cells.Item(92, "H").Value = 2700.9
$ws.Cells.Item(92, "I").Value = 2913.625
$ws.Cells.Item(92, "K").Value = 2913.625
$ws.Cells.Item(92, "M").Value = -1665.625

# Row 111 (Leve Item ID 27768)
$ws.Cells.Item(111, "H").Value = 1334.5
$ws.Cells.Item(111, "I").Value = 619.75
$ws.Cells.Item(111, "J").Value = 2764
$ws.Cells.Item(111, "K").Value = 1859.25
$ws.Cells.Item(111, "L").Value = 8292
$ws.Cells.Item(111, "M").Value = 1207.75
$ws.Cells.Item(111, "N").Value = -14426

# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, "H").Value = 6733.4443
$ws.Cells.Item(132, "I").Value = 6699.375
$ws.Cells.Item(132, "J").Value = 7006
$ws.Cells.Item(132, "K").Value = 20098.125
$ws.Cells.Item(132, "L").Value = 21018
$ws.Cells.Item(132, "M").Value = -17568.125
$ws.Cells.Item(132, "N").Value = -26078

# Row 137 (Leve Item ID 44013)
$ws.Cells.Item(137, "H").Value = 1079.175
$ws.Cells.Item(137, "I").Value = 878.5625
$ws.Cells.Item(137, "J").Value = 1881.625
$ws.Cells.Item(137, "K").Value = 2635.6875
$ws.Cells.Item(137, "L").Value = 5644.875
$ws.Cells.Item(137, "M").Value = -85.6875
$ws.Cells.Item(137, "N").Value = -10744.875

# Row 138 (Leve Item ID 44169)
$ws.Cells.Item(138, "H").Value = 4106.3
$ws.Cells.Item(138, "I").Value = 858.75
$ws.Cells.Item(138, "J").Value = 4918.1875
$ws.Cells.Item(138, "K").Value = 2576.25
$ws.Cells.Item(138, "L").Value = 14754.5625
$ws.Cells.Item(138, "M").Value = 2563.75
$ws.Cells.Item(138, "N").Value = -25034.5625

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Cells.Item(2, "H").Value = 45314.348
$ws.Cells.Item(2, "I").Value = 64377.438
$ws.Cells.Item(2, "J").Value = 1741.5714
$ws.Cells.Item(2, "K").Value = 64377.438
$ws.Cells.Item(2, "L").Value = 1741.5714
$ws.Cells.Item(2, "M").Value = -64264.438
$ws.Cells.Item(2, "N").Value = -1967.5714

# Row 45 (Leve Item ID 27714)
$ws.Cells.Item(45, "H").Value = 1704.4
$ws.Cells.Item(45, "I").Value = 880.5
$ws.Cells.Item(45, "J").Value = 5000
$ws.Cells.Item(45, "K").Value = 880.5
$ws.Cells.Item(45, "L").Value = 5000
$ws.Cells.Item(45, "M").Value = -503.5
$ws.Cells.Item(45, "N").Value = -5754

# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, "H").Value = 1116.6923
$ws.Cells.Item(74, "I").Value = 1143.0834
$ws.Cells.Item(74, "K").Value = 1143.0834
$ws.Cells.Item(74, "M").Value = -269.0834

# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, "H").Value = 1116.6923
$ws.Cells.Item(77, "I").Value = 1143.0834
$ws.Cells.Item(77, "K").Value = 5715.416999999999
$ws.Cells.Item(77, "M").Value = -1347.416999999999

# Row 116 (Leve Item ID 27713)
$ws.Cells.Item(116, "H").Value = 45314.348
$ws.Cells.Item(116, "I").Value = 64377.438
$ws.Cells.Item(116, "J").Value = 1741.5714
$ws.Cells.Item(116, "K").Value = 64377.438
$ws.Cells.Item(116, "L").Value = 1741.5714
$ws.Cells.Item(116, "M").Value = -62083.438
$ws.Cells.Item(116, "N").Value = -6329.5714

# Row 123 (Leve Item ID 34107)
$ws.Cells.Item(123, "H").Value = 0
$ws.Cells.Item(123, "J").Value = 0
$ws.Cells.Item(123, "L").Value = 0
$ws.Cells.Item(123, "N").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Cells.Item(3, "H").Value = 45314.348
$ws.Cells.Item(3, "I").Value = 64377.438
$ws.Cells.Item(3, "J").Value = 1741.5714
$ws.Cells.Item(3, "K").Value = 64377.438
$ws.Cells.Item(3, "L").Value = 1741.5714
$ws.Cells.Item(3, "M").Value = -64263.438
$ws.Cells.Item(3, "N").Value = -1969.5714

# Row 99 (Leve Item ID 19943)
$ws.Cells.Item(99, "H").Value = 50001508
$ws.Cells.Item(99, "I").Value = 55557130
$ws.Cells.Item(99, "J").Value = 940
$ws.Cells.Item(99, "K").Value = 55557130
$ws.Cells.Item(99, "L").Value = 940
$ws.Cells.Item(99, "M").Value = -55555632
$ws.Cells.Item(99, "N").Value = -3936

# Row 107 (Leve Item ID 27706)
$ws.Cells.Item(107, "H").Value = 15232.7
$ws.Cells.Item(107, "I").Value = 1815.875
$ws.Cells.Item(107, "J").Value = 68900
$ws.Cells.Item(107, "K").Value = 1815.875
$ws.Cells.Item(107, "L").Value = 68900
$ws.Cells.Item(107, "M").Value = 104.125
$ws.Cells.Item(107, "N").Value = -72740

# Row 109 (Leve Item ID 27096)
$ws.Cells.Item(109, "H").Value = 0
$ws.Cells.Item(109, "J").Value = 0
$ws.Cells.Item(109, "L").Value = 0
$ws.Cells.Item(109, "N").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, "H").Value = 2716.359
$ws.Cells.Item(31, "I").Value = 2949.1724
$ws.Cells.Item(31, "J").Value = 2041.2
$ws.Cells.Item(31, "K").Value = 2949.1724
$ws.Cells.Item(31, "L").Value = 2041.2
$ws.Cells.Item(31, "M").Value = -2654.1724
$ws.Cells.Item(31, "N").Value = -2631.2

# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, "H").Value = 2716.359
$ws.Cells.Item(34, "I").Value = 2949.1724
$ws.Cells.Item(34, "J").Value = 2041.2
$ws.Cells.Item(34, "K").Value = 2949.1724
$ws.Cells.Item(34, "L").Value = 2041.2
$ws.Cells.Item(34, "M").Value = -2747.1724
$ws.Cells.Item(34, "N").Value = -2445.2

# Row 107 (Leve Item ID 27689)
$ws.Cells.Item(107, "H").Value = 419.0435
$ws.Cells.Item(107, "I").Value = 595
$ws.Cells.Item(107, "J").Value = 402.2857
$ws.Cells.Item(107, "K").Value = 595
$ws.Cells.Item(107, "L").Value = 402.2857
$ws.Cells.Item(107, "M").Value = 1325
$ws.Cells.Item(107, "N").Value = -4242.2857

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Cells.Item(12, "H").Value = 62.333332
$ws.Cells.Item(12, "J").Value = 68
$ws.Cells.Item(12, "L").Value = 204
$ws.Cells.Item(12, "N").Value = -550

# Row 92 (Leve Item ID 19841)
$ws.Cells.Item(92, "H").Value = 518.75
$ws.Cells.Item(92, "I").Value = 300
$ws.Cells.Item(92, "J").Value = 591.6667
$ws.Cells.Item(92, "K").Value = 900
$ws.Cells.Item(92, "L").Value = 1775.0001
$ws.Cells.Item(92, "M").Value = 348
$ws.Cells.Item(92, "N").Value = -4271.0001

# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, "H").Value = 13424.566
$ws.Cells.Item(131, "I").Value = 56407.777
$ws.Cells.Item(131, "J").Value = 1521.5231
$ws.Cells.Item(131, "K").Value = 169223.331
$ws.Cells.Item(131, "L").Value = 4564.5693
$ws.Cells.Item(131, "M").Value = -164183.331
$ws.Cells.Item(131, "N").Value = -14644.5693

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Cells.Item(70, "H").Value = 5831.3335
$ws.Cells.Item(70, "I").Value = 4996
$ws.Cells.Item(70, "J").Value = 6666.6665
$ws.Cells.Item(70, "K").Value = 4996
$ws.Cells.Item(70, "L").Value = 6666.6665
$ws.Cells.Item(70, "M").Value = -4726
$ws.Cells.Item(70, "N").Value = -7206.6665

# Row 73 (Leve Item ID 14146)
$ws.Cells.Item(73, "H").Value = 5831.3335
$ws.Cells.Item(73, "I").Value = 4996
$ws.Cells.Item(73, "J").Value = 6666.6665
$ws.Cells.Item(73, "K").Value = 4996
$ws.Cells.Item(73, "L").Value = 6666.6665
$ws.Cells.Item(73, "M").Value = -4060
$ws.Cells.Item(73, "N").Value = -8538.666499999999

# Row 113 (Leve Item ID 27710)
$ws.Cells.Item(113, "H").Value = 4991.5835
$ws.Cells.Item(113, "I").Value = 5892.1113
$ws.Cells.Item(113, "J").Value = 2290
$ws.Cells.Item(113, "K").Value = 5892.1113
$ws.Cells.Item(113, "L").Value = 2290
$ws.Cells.Item(113, "M").Value = -3722.1113
$ws.Cells.Item(113, "N").Value = -6630

# Row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, "H").Value = 1881001.2
$ws.Cells.Item(122, "I").Value = 2632840.2
$ws.Cells.Item(122, "K").Value = 7898520.600000001
$ws.Cells.Item(122, "M").Value = -7896070.600000001

# Row 123 (Leve Item ID 34150)
$ws.Cells.Item(123, "H").Value = 10763.8
$ws.Cells.Item(123, "J").Value = 10763.8
$ws.Cells.Item(123, "L").Value = 10763.8
$ws.Cells.Item(123, "N").Value = -15663.8

# Row 133 (Leve Item ID 41854)
$ws.Cells.Item(133, "H").Value = 44000
$ws.Cells.Item(133, "J").Value = 44000
$ws.Cells.Item(133, "L").Value = 44000
$ws.Cells.Item(133, "N").Value = -54120

$ws = $wb.Worksheets.Item("LTW")
# Row 29 (Leve Item ID 3576)
$ws.Cells.Item(29, "H").Value = 0
$ws.Cells.Item(29, "J").Value = 0
$ws.Cells.Item(29, "L").Value = 0
$ws.Cells.Item(29, "N").ClearContents()

# Row 55 (Leve Item ID 5284)
$ws.Cells.Item(55, "H").Value = 500
$ws.Cells.Item(55, "I").Value = 400
$ws.Cells.Item(55, "J").Value = 550
$ws.Cells.Item(55, "K").Value = 400
$ws.Cells.Item(55, "L").Value = 550
$ws.Cells.Item(55, "M").Value = -227
$ws.Cells.Item(55, "N").Value = -896

# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, "H").Value = 7187.2915
$ws.Cells.Item(122, "I").Value = 7599.75
$ws.Cells.Item(122, "K").Value = 22799.25
$ws.Cells.Item(122, "M").Value = -20349.25

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, "H").Value = 1641.7258
$ws.Cells.Item(132, "I").Value = 1407.2885
$ws.Cells.Item(132, "J").Value = 2860.8
$ws.Cells.Item(132, "K").Value = 4221.8655
$ws.Cells.Item(132, "L").Value = 8582.400000000001
$ws.Cells.Item(132, "M").Value = -1691.8655
$ws.Cells.Item(132, "N").Value = -13642.4

# Row 136 (Leve Item ID 44060)
$ws.Cells.Item(136, "H").Value = 4230.4375
$ws.Cells.Item(136, "I").Value = 2013.3572
$ws.Cells.Item(136, "J").Value = 19750
$ws.Cells.Item(136, "K").Value = 6040.071599999999
$ws.Cells.Item(136, "L").Value = 59250
$ws.Cells.Item(136, "M").Value = -3490.071599999999
$ws.Cells.Item(136, "N").Value = -64350

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Cells.Item(107, "H").Value = 644.26086
$ws.Cells.Item(107, "I").Value = 643.375
$ws.Cells.Item(107, "K").Value = 1930.125
$ws.Cells.Item(107, "M").Value = -10.125

# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, "H").Value = 3555.6365
$ws.Cells.Item(132, "I").Value = 3626
$ws.Cells.Item(132, "J").Value = 3626
$ws.Cells.Item(132, "K").Value = 10878
$ws.Cells.Item(132, "L").Value = 10546.2855
$ws.Cells.Item(132, "M").Value = -8348
$ws.Cells.Item(132, "N").Value = -15606.2855
